# Grading pass: "kalyankar to pusapati done"
#
# The TA fills in the "Total Points" column (E) for the "Generic" class
# section (rows 3-6) and the "Customer Class" section (rows 10-14),
# awarding full marks by mirroring the "Points for grading" column (D).
# The section/grand totals (E7, E15, E38) are formulas and recompute
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the cursor/selection to E15 (the running total for the section just
# graded), matching where the grader's view was left.
$ws.Range("E15").Select()
